# Update the NATMI LR-pair TPM values (Wnt5a-Fzd3 sheet) with the new TPM
# derived numbers, per commit "update scripts wuth new tpm".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.02354566666666667
$ws.Range("H2").Value = 0.07063700000000001
$ws.Range("I2").Value = 0.002815555392485919
$ws.Range("J2").Value = 0.002815555392485918
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.196431
$ws.Range("N2").Value = 0.589293
$ws.Range("O2").Value = 0.09717285149889213
$ws.Range("P2").Value = 0.09717285149889213
$ws.Range("Q2").Value = 0.004625098849000001
$ws.Range("R2").Value = 0.041625889641
$ws.Range("S2").Value = 0.0002735955460409391
$ws.Range("T2").Value = 0.0002735955460409391

$ws.Range("G3").Value = 0.02354566666666667
$ws.Range("H3").Value = 0.07063700000000001
$ws.Range("I3").Value = 0.002815555392485919
$ws.Range("J3").Value = 0.002815555392485918
$ws.Range("M3").Value = 0.4307096666666667
$ws.Range("O3").Value = 0.2130686423127578
$ws.Range("P3").Value = 0.2130686423127578
$ws.Range("Q3").Value = 0.01014134624144445
$ws.Range("R3").Value = 0.09127211617300002
$ws.Range("S3").Value = 0.0005999065648333387
$ws.Range("T3").Value = 0.0005999065648333386

$ws.Range("G4").Value = 0.02354566666666667
$ws.Range("H4").Value = 0.07063700000000001
$ws.Range("I4").Value = 0.002815555392485919
$ws.Range("J4").Value = 0.002815555392485918
$ws.Range("O4").Value = 0.68975850618835
$ws.Range("P4").Value = 0.68975850618835
$ws.Range("Q4").Value = 0.032830170401
$ws.Range("R4").Value = 0.295471533609
$ws.Range("S4").Value = 0.001942053281611641
$ws.Range("T4").Value = 0.00194205328161164

$ws.Range("I5").Value = 0.9868456480383168
$ws.Range("J5").Value = 0.9868456480383166
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.196431
$ws.Range("N5").Value = 0.589293
$ws.Range("O5").Value = 0.09717285149889213
$ws.Range("P5").Value = 0.09717285149889213
$ws.Range("Q5").Value = 1.621086441085
$ws.Range("R5").Value = 14.589777969765
$ws.Range("S5").Value = 0.09589460560915532
$ws.Range("T5").Value = 0.0958946056091553

$ws.Range("I6").Value = 0.9868456480383168
$ws.Range("J6").Value = 0.9868456480383166
$ws.Range("M6").Value = 0.4307096666666667
$ws.Range("O6").Value = 0.2130686423127578
$ws.Range("P6").Value = 0.2130686423127578
$ws.Range("Q6").Value = 3.554518383949445
$ws.Range("S6").Value = 0.2102658623997778
$ws.Range("T6").Value = 0.2102658623997778

$ws.Range("I7").Value = 0.9868456480383168
$ws.Range("J7").Value = 0.9868456480383166
$ws.Range("O7").Value = 0.68975850618835
$ws.Range("P7").Value = 0.68975850618835
$ws.Range("S7").Value = 0.6806851800293836
$ws.Range("T7").Value = 0.6806851800293835

$ws.Range("I8").Value = 0.0103387965691973
$ws.Range("J8").Value = 0.0103387965691973
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.196431
$ws.Range("N8").Value = 0.589293
$ws.Range("O8").Value = 0.09717285149889213
$ws.Range("P8").Value = 0.09717285149889213
$ws.Range("Q8").Value = 0.016983489737
$ws.Range("R8").Value = 0.152851407633
$ws.Range("S8").Value = 0.001004650343695865
$ws.Range("T8").Value = 0.001004650343695865

$ws.Range("I9").Value = 0.0103387965691973
$ws.Range("J9").Value = 0.0103387965691973
$ws.Range("M9").Value = 0.4307096666666667
$ws.Range("O9").Value = 0.2130686423127578
$ws.Range("P9").Value = 0.2130686423127578
$ws.Range("R9").Value = 0.335153712149
$ws.Range("S9").Value = 0.002202873348146668
$ws.Range("T9").Value = 0.002202873348146668

$ws.Range("I10").Value = 0.0103387965691973
$ws.Range("J10").Value = 0.0103387965691973
$ws.Range("O10").Value = 0.68975850618835
$ws.Range("P10").Value = 0.68975850618835
$ws.Range("S10").Value = 0.00713127287735477
$ws.Range("T10").Value = 0.00713127287735477
